$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the Date value (row 8, column B) to the new timestamp.
$ws.Range("B8").Value = "2025-07-17T14:35:50+00:00"

# Fill in the previously empty Title value (row 5, column B)
# with the same text already used for Name (row 4, column B).
$ws.Range("B5").Value = "CapaciteSavoirfaire"
